$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update milestone text in A14 ("TDD check " -> "Demonstration – All tests should pass")
$ws.Range("A14").Value = "Demonstration" + [char]0x0020 + [char]0x2013 + " All tests should pass"

# Move the active selection from C17 to A14
$null = $ws.Range("A14").Select()
